$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.136.42'
$ws.Range('E2').Value = '  -1.70%  '

$ws.Range('D3').Value = '1.806.17'
$ws.Range('E3').Value = '  -1.74%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = '309.72'
$ws.Range('E5').Value = '  -1.63%  '

$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('D7').Value = '0.4235'
$ws.Range('E7').Value = '  -1.16%  '

$ws.Range('D8').Value = '0.3627'
$ws.Range('E8').Value = '  -1.48%  '

$ws.Range('D9').Value = '0.07222'
$ws.Range('E9').Value = '  -0.97%  '

$ws.Range('D10').Value = '0.8481'
$ws.Range('E10').Value = '  -2.76%  '

$ws.Range('D11').Value = '20.33'
$ws.Range('E11').Value = '  -2.17%  '

$ws.Range('D12').Value = '1.816.87'
$ws.Range('E12').Value = '  -4.09%  '

$ws.Range('D13').Value = '5.321'
$ws.Range('E13').Value = '  -2.13%  '

$ws.Range('D14').Value = '6.413'
$ws.Range('E14').Value = '  -2.29%  '

$ws.Range('D15').Value = '0.06819'
$ws.Range('E15').Value = '  -1.75%  '

$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '81.50'
$ws.Range('E16').Value = '  +1.46%  '

$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  +0.33%  '

$ws.Range('D18').Value = '0.000008763'
$ws.Range('E18').Value = '  -2.43%  '

$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  -0.02%  '

$ws.Range('D20').Value = '15.04'
$ws.Range('E20').Value = '  -3.05%  '

$ws.Range('D21').Value = '27.010.65'
$ws.Range('E21').Value = '  -2.88%  '

$ws.Range('D22').Value = '5.109'
$ws.Range('E22').Value = '  -1.01%  '

$ws.Range('E23').Value = '  +2.09%  '

$ws.Range('D24').Value = '1.998.63'
$ws.Range('E24').Value = '  -7.26%  '

$ws.Range('D25').Value = '1.952'
$ws.Range('E25').Value = '  -1.63%  '

$ws.Range('D26').Value = '153.23'
$ws.Range('E26').Value = '  -0.21%  '

$ws.Range('D27').Value = '18.32'
$ws.Range('E27').Value = '  -2.36%  '

$ws.Range('D28').Value = '5.051'
$ws.Range('E28').Value = '  -3.70%  '

$ws.Range('D29').Value = '115.06'
$ws.Range('E29').Value = '  +0.55%  '

$ws.Range('D30').Value = '1.657'
$ws.Range('E30').Value = '  -10.77%  '

$ws.Range('D31').Value = '0.08958'
$ws.Range('E31').Value = '  +0.91%  '

$ws.Range('D32').Value = '0.7378'
$ws.Range('E32').Value = '  -4.35%  '

$ws.Range('D33').Value = '2.882'
$ws.Range('E33').Value = '  -2.65%  '

$ws.Range('D34').Value = '4.382'
$ws.Range('E34').Value = '  -3.94%  '

$ws.Range('D35').Value = '1.101'
$ws.Range('E35').Value = '  -4.24%  '

$ws.Range('D36').Value = '1.003'
$ws.Range('E36').Value = '  +0.09%  '

$ws.Range('D37').Value = '1.084'
$ws.Range('E37').Value = '  -1.22%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01912'
$ws.Range('E38').Value = '  -1.97%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05154'
$ws.Range('E39').Value = '  -4.10%  '

$ws.Range('D40').Value = '0.1633'
$ws.Range('E40').Value = '  -2.58%  '

$ws.Range('D41').Value = '0.4982'
$ws.Range('E41').Value = '  -2.44%  '

$ws.Range('D42').Value = '2.617'
$ws.Range('E42').Value = '  -7.34%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '6.059'
$ws.Range('E43').Value = '  -9.10%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '8.114'
$ws.Range('E44').Value = '  -4.71%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '10.29'
$ws.Range('E45').Value = '  -2.01%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '104.90'
$ws.Range('E46').Value = '  -1.27%  '

$ws.Range('D47').Value = '1.003'
$ws.Range('E47').Value = '  +0.13%  '

$ws.Range('D48').Value = '0.06319'
$ws.Range('E48').Value = '  -3.21%  '

$ws.Range('D49').Value = '0.4543'
$ws.Range('E49').Value = '  -3.91%  '

$ws.Range('D50').Value = '1.604'
$ws.Range('E50').Value = '  -1.53%  '

$ws.Range('D51').Value = '1.748'
$ws.Range('E51').Value = '  -1.96%  '
